# Applies the "cryptos list" refresh captured in the commit diff:
# updated Price (D) / Volume-1h (E) figures for every row, plus a handful of
# coins that changed rank position (so Coin name / Link in B & C moved too).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.104.90'
$ws.Range('E2').Value = '  -1.94%  '

# Row 3
$ws.Range('D3').Value = '1.834.24'
$ws.Range('E3').Value = '  -3.35%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.24%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.19'
$ws.Range('E5').Value = '  -4.64%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  +0.28%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4590'
$ws.Range('E7').Value = '  -4.94%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2660'
$ws.Range('E8').Value = '  -6.79%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06199'
$ws.Range('E9').Value = '  -5.46%  '

# Row 10
$ws.Range('D10').Value = '1.846.33'
$ws.Range('E10').Value = '  -4.32%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07333'
$ws.Range('E11').Value = '  -1.75%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.85'
$ws.Range('E12').Value = '  -5.39%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.841'
$ws.Range('E13').Value = '  -5.29%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '82.71'
$ws.Range('E14').Value = '  -6.08%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6133'
$ws.Range('E15').Value = '  -8.13%  '

# Row 16
$ws.Range('D16').Value = '30.020.72'
$ws.Range('E16').Value = '  -2.17%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.003'
$ws.Range('E17').Value = '  +0.20%  '

# Row 18
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '224.48'
$ws.Range('E18').Value = '  -3.00%  '

# Row 19
$ws.Range('B19').Value = 'BinanceUSD'
$ws.Range('C19').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.005'
$ws.Range('E19').Value = '  +0.43%  '

# Row 20
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.29'
$ws.Range('E20').Value = '  -7.82%  '

# Row 21
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000007177'
$ws.Range('E21').Value = '  -5.72%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.832'
$ws.Range('E22').Value = '  -8.67%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.819'
$ws.Range('E23').Value = '  -6.83%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '164.38'
$ws.Range('E24').Value = '  -3.26%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.047'
$ws.Range('E25').Value = '  -3.17%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.39'
$ws.Range('E26').Value = '  -7.39%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.838'
$ws.Range('E27').Value = '  -6.62%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1009'
$ws.Range('E28').Value = '  -1.06%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.376'
$ws.Range('E29').Value = '  -2.03%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.021'
$ws.Range('E30').Value = '  -7.71%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.733'
$ws.Range('E31').Value = '  -7.42%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.04790'
$ws.Range('E32').Value = '  -6.41%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.125'
$ws.Range('E33').Value = '  -7.64%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.6875'
$ws.Range('E34').Value = '  -9.50%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.002'
$ws.Range('E35').Value = '  -0.18%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.708'
$ws.Range('E36').Value = '  -0.01%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01805'
$ws.Range('E37').Value = '  -4.27%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.600'
$ws.Range('E38').Value = '  -2.30%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8893'
$ws.Range('E39').Value = '  -3.51%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.909'
$ws.Range('E40').Value = '  -8.26%  '

# Row 41
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9973'
$ws.Range('E41').Value = '  -0.71%  '

# Row 42
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '102.95'
$ws.Range('E42').Value = '  -3.89%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.425'
$ws.Range('E43').Value = '  -5.28%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3966'
$ws.Range('E44').Value = '  -7.93%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.843'
$ws.Range('E45').Value = '  -8.08%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1169'
$ws.Range('E46').Value = '  -8.41%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '58.70'
$ws.Range('E47').Value = '  -9.21%  '

# Row 48
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05547'
$ws.Range('E48').Value = '  -2.37%  '

# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.356'
$ws.Range('E49').Value = '  -6.68%  '

# Row 50
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '32.29'
$ws.Range('E50').Value = '  -4.83%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.351'
$ws.Range('E51').Value = '  -9.63%  '
